$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns at H:I (Volume / Volume Uncertainty).
#    This shifts the old H,I,J,K (Density, Density Unc, Elemental Purity,
#    Note) to J,K,L,M and fixes up every formula reference automatically.
# ---------------------------------------------------------------------------
$ws.Range("H1:I1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. Re-label column B (was "Circumference (mm)", now "Diameter (mm)") and
#    give the whole header row its new wrapped / centered style + height.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value2 = "Diameter (mm)"

$ws.Range("H1").Value2 = "Volume [cm^3]"
$ws.Range("I1").Value2 = "Volume Uncertainty"
$ws.Range("K1").Value2 = "Density Uncertainty"

# Build the wrapped/centered style once on B1, then fan it out to the rest
# of the header row via a format-only paste so the stylesheet ends up with
# a single new cellXfs entry (not one per property assignment).
$b1 = $ws.Range("B1")
$b1.HorizontalAlignment = -4108   # xlCenter
$b1.VerticalAlignment = -4108     # xlCenter
$b1.WrapText = $true
$b1.Copy()
$ws.Range("C1:M1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows(1).RowHeight = 30

# ---------------------------------------------------------------------------
# 3. New Volume + Volume Uncertainty formulas. Row 2 is entered on its own
#    (standalone formula, matching the pre-existing pattern in this sheet
#    where the first data row is typed individually), rows 3-6 are filled
#    as one range so the engine groups them into a shared formula.
# ---------------------------------------------------------------------------
$ws.Range("H2").Formula = "=PI()*B2^2/4*D2/1000"
$ws.Range("I2").Formula = "=SQRT((C2/B2)^2+(E2/D2)^2)*H2"

$ws.Range("H3:H6").Formula = "=PI()*B3^2/4*D3/1000"
$ws.Range("I3:I6").Formula = "=SQRT((C3/B3)^2+(E3/D3)^2)*H3"

# ---------------------------------------------------------------------------
# 4. Density Uncertainty (old column I, now K) needs to reference the
#    (also shifted) Density column J instead of H - re-enter it as a shared
#    formula across the same K3:K6 run it had before the insert.
# ---------------------------------------------------------------------------
$ws.Range("K3:K6").Formula = "=SQRT((C3/B3)^2+(E3/D3)^2+(G3/F3)^2)*J3"

# ---------------------------------------------------------------------------
# 5. Extra note cell that now lives at K15.
# ---------------------------------------------------------------------------
$ws.Range("K15").Value2 = " "

# ---------------------------------------------------------------------------
# 6. Column widths - best effort match of the target layout (the engine
#    quantizes ColumnWidth to 1/6-character steps, so these are the closest
#    achievable inputs to the recorded target widths).
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 13.333333333333334
$ws.Columns("C").ColumnWidth = 11.0
$ws.Columns("D").ColumnWidth = 8.666666666666666
$ws.Columns("E").ColumnWidth = 10.833333333333334
$ws.Columns("F").ColumnWidth = 9.5
$ws.Columns("G").ColumnWidth = 12.0
$ws.Columns("H").ColumnWidth = 9.666666666666666
$ws.Columns("I").ColumnWidth = 11.666666666666666
$ws.Columns("J").ColumnWidth = 11.333333333333334
$ws.Columns("K").ColumnWidth = 13.5
$ws.Columns("L").ColumnWidth = 13.5
$ws.Columns("M").ColumnWidth = 11.666666666666666

# ---------------------------------------------------------------------------
# 7. Misc view / print state.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait

$ws.Range("I2").Select()
